$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(9.623814262023211,0.08,34.84),
  @(2.676119402985066,0,11.46),
  @(3.656451077943601,0,14.61333333333333),
  @(2.947197346600325,0.006666666666666667,11.69333333333333),
  @(6.216285240464338,0.08666666666666667,22.79333333333333),
  @(9.606733001658364,0.1066666666666667,35.59333333333333),
  @(12.56046434494194,0.2400000000000001,47.50666666666668),
  @(5.828723051409612,0.006666666666666667,21.18),
  @(12.45243781094526,0.2066666666666667,45.74666666666668),
  @(11.64092868988391,0.2600000000000001,41.57333333333334),
  @(8.601459369817571,0.3133333333333335,28.79333333333334),
  @(12.66971807628523,0.2733333333333334,45.46666666666665),
  @(11.78885572139303,0.2600000000000001,43.87333333333333),
  @(3.985737976782742,0,15.58666666666666),
  @(6.569485903814253,0.02666666666666667,24.42666666666667),
  @(7.831575456053065,0.09333333333333335,27.52),
  @(3.385140961857366,0,15.08666666666667),
  @(11.66222222222222,0.3866666666666667,40.96666666666666),
  @(6.69339966832503,0.08,24.34),
  @(4.80487562189054,0.01333333333333333,17.54666666666666),
  @(12.93724709784411,0.2200000000000001,48.4866666666667),
  @(3.745870646766159,0.006666666666666667,16.33333333333333),
  @(6.797910447761182,0.05333333333333333,24.46666666666667),
  @(5.44169154228855,0.05999999999999999,20.25999999999999),
  @(4.779900497512427,0.02,18.98666666666666),
  @(12.7678606965174,0.3400000000000001,45.22),
  @(9.268888888888874,0.1800000000000001,34.88666666666666),
  @(4.20590381426201,0,17.40666666666667),
  @(8.089883913764499,0.2600000000000001,29.85333333333333),
  @(12.07824212271973,0.1,44.70666666666667),
  @(11.56673300165837,0.2200000000000001,43.94666666666668),
  @(11.78822553897181,0.2800000000000001,43.61999999999999),
  @(11.38126036484244,0.3266666666666668,42.02),
  @(6.83830845771143,0.09333333333333334,24.75999999999999),
  @(9.410248756218895,0.1266666666666667,34.84666666666667),
  @(4.692504145936973,0.02666666666666667,18.59333333333333),
  @(11.37698175787728,0.1533333333333334,38.96666666666666),
  @(9.134726368159191,0.2133333333333334,33.17333333333332),
  @(3.59728026533995,0,13.86666666666666),
  @(8.922354892205625,0.08666666666666667,33.84666666666666),
  @(11.65273631840795,0.3200000000000001,44.03333333333333),
  @(8.278872305140951,0.2466666666666668,30.31333333333333),
  @(10.18162520729684,0.2066666666666667,39.10000000000001),
  @(11.21933665008291,0.1200000000000001,42.05333333333333),
  @(11.83323383084576,0.1066666666666667,46.16666666666666),
  @(12.90567164179104,0.2733333333333334,45.25333333333333),
  @(8.551608623548908,0.1466666666666666,31.95333333333333),
  @(3.789751243781078,0.006666666666666667,14.13333333333334),
  @(5.205538971807621,0.006666666666666667,19.17333333333334),
  @(2.730878938640122,0.006666666666666667,11.54),
  @(11.35333333333333,0.2,42.34000000000001),
  @(8.006368159203971,0.04666666666666667,29.48666666666667),
  @(3.719303482587054,0.006666666666666667,14.83999999999999),
  @(4.725671641791034,0.04666666666666666,18.49333333333333),
  @(11.33532338308457,0.2200000000000001,43.34),
  @(11.8751575456053,0.1933333333333334,44.68666666666668),
  @(4.907396351575446,0.06666666666666667,19.74666666666667),
  @(11.27333333333332,0.1933333333333334,39.16666666666666),
  @(4.288026533996672,0,17.66666666666666),
  @(12.45615257048092,0.3600000000000001,46.17999999999999),
  @(11.59535655058042,0.2066666666666667,42.83333333333334),
  @(13.16464344941957,0.3800000000000002,51.25999999999998),
  @(3.981393034825859,0,16.15999999999999),
  @(12.56391376451077,0.38,45.76000000000001),
  @(10.54248756218904,0.3266666666666667,36.88666666666666),
  @(10.41860696517412,0.2266666666666667,40.23999999999999),
  @(4.57014925373133,0.006666666666666667,17.61333333333333),
  @(12.29532338308457,0.2533333333333335,46.96666666666668),
  @(3.678905472636804,0.01333333333333333,14.06666666666666),
  @(7.322122719734653,0.08666666666666667,28.07333333333333),
  @(7.886368159203971,0.1,28.74666666666667)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
